$d = $word.ActiveDocument

function Replace-Text($search, $replace) {
    $rng = $d.Content
    $rng.Find.Execute($search, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

function Color-Text($search, $colorMode) {
    # colorMode: "black" -> w:color val=000000 ; "none" -> no w:color element
    $rng = $d.Content
    $found = $rng.Find.Execute($search, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $piece = $d.Range($rng.Start, $rng.End)
        if ($colorMode -eq "black") {
            $piece.Font.Color = 0
        } else {
            $piece.Font.ColorIndex = 0
        }
    }
}

# --- Block 1: ", keep the space for shadow" -> ", keep the space of shadow" ---
Replace-Text ", keep the space for shadow" ", keep the space of shadow"
Color-Text "of shadow" "black"
Color-Text "of " "none"

# --- Block 2: "alone by themselves &" -> "for these alone, &" ---
Replace-Text "alone by themselves &" "for these alone, &"
Color-Text "for these alone, &" "black"
Color-Text "for " "none"
Color-Text "se alone" "none"

# --- Block 3: " for the lights &" -> " also individually that of lights &" ---
Replace-Text " for the lights &" " also individually that of lights &"
Color-Text " also individually that of lights &" "black"
Color-Text "individually that" "none"
Color-Text "that of lights" "none"

# --- Block 4: " highlights individually as well, " -> " highlights, " ---
Replace-Text " highlights individually as well, " " highlights, "

# --- Block 5: " a colo" -> " one colo" ---
Replace-Text " a colo" " one colo"
Color-Text " one colo" "black"
Color-Text "one" "none"

# --- Block 6: "economise" -> "economize" ---
Replace-Text "economise" "economize"

# --- Block 7: "neatly, which " -> "neatly. Which " ---
Replace-Text "neatly, which " "neatly. Which "
Color-Text "neatly. Which " "black"
Color-Text "Which " "none"

# --- Block 8: ". Use, then, the " -> ". Therefore use the " ---
Replace-Text ". Use, then, the " ". Therefore use the "
